# Update countries & provincias Spain
# Applies the 20-May-2020 10:35 data refresh to the "Pais" sheet:
#  - updates the "last updated" timestamp string
#  - refreshes case counters for several countries
#  - re-sorts the three tied/re-ranked country pairs (Corea del Sur / Dinamarca,
#    El Salvador / Lituania, Nueva Caledonia / Santa Lucia / Belice) by moving
#    their case-count rows so the table stays ordered by total cases desc.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header timestamp -------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 20 de Mayo de 2020 a las 10:35"

# --- Row 34: Polonia (rank unchanged) ----------------------------------
$ws.Range("B34").Value = 19569
$ws.Range("C34").Value = 301
$ws.Range("E34").Value = 10433
$ws.Range("G34").Value = 5
$ws.Range("H34").Value = 953

# --- Row 36: Indonesia (rank unchanged) --------------------------------
$ws.Range("B36").Value = 19189
$ws.Range("C36").Value = 693
$ws.Range("D36").Value = 4575
$ws.Range("E36").Value = 13372
$ws.Range("G36").Value = 21
$ws.Range("H36").Value = 1242

# --- Row 46: Filipinas (rank unchanged) --------------------------------
$ws.Range("B46").Value = 13221
$ws.Range("C46").Value = 279
$ws.Range("D46").Value = 2932
$ws.Range("E46").Value = 9447
$ws.Range("G46").Value = 5
$ws.Range("H46").Value = 842

# --- Rows 47/48: Dinamarca overtakes Corea del Sur ---------------------
# Dinamarca's refreshed totals (11117) now exceed Corea del Sur's (11110),
# so Dinamarca takes row 47 and Corea del Sur (unchanged) drops to row 48.
$ws.Range("A47").Value = "Dinamarca"
$ws.Range("B47").Value = 11117
$ws.Range("C47").Value = 73
$ws.Range("D47").Value = 9416
$ws.Range("E47").Value = 1150
$ws.Range("F47").Value = 0
$ws.Range("G47").Value = 0
$ws.Range("H47").Value = 551

$ws.Range("A48").Value = "Corea del Sur"
$ws.Range("B48").Value = 11110
$ws.Range("C48").Value = 32
$ws.Range("D48").Value = 10066
$ws.Range("E48").Value = 781
$ws.Range("F48").Value = 0
$ws.Range("G48").Value = 0
$ws.Range("H48").Value = 263

# --- Row 89: Estonia (rank unchanged) ----------------------------------
$ws.Range("B89").Value = 1794
$ws.Range("C89").Value = 3
$ws.Range("D89").Value = 945
$ws.Range("E89").Value = 785

# --- Rows 92/93: Lituania overtakes El Salvador ------------------------
# Lituania's refreshed totals (1577) now exceed El Salvador's (1571),
# so Lituania takes row 92 and El Salvador (unchanged) drops to row 93.
$ws.Range("A92").Value = "Lituania"
$ws.Range("B92").Value = 1577
$ws.Range("C92").Value = 15
$ws.Range("D92").Value = 1049
$ws.Range("E92").Value = 468
$ws.Range("F92").Value = 0
$ws.Range("G92").Value = 0
$ws.Range("H92").Value = 60

$ws.Range("A93").Value = "El Salvador"
$ws.Range("B93").Value = 1571
$ws.Range("C93").Value = 73
$ws.Range("D93").Value = 531
$ws.Range("E93").Value = 1009
$ws.Range("F93").Value = 0
$ws.Range("G93").Value = 1
$ws.Range("H93").Value = 31

# --- Row 159: Mozambique (rank unchanged) ------------------------------
$ws.Range("D159").Value = 48
$ws.Range("E159").Value = 98

# --- Rows 195/196/197: Belice jumps ahead of Nueva Caledonia & Santa Lucia
# All three are tied on total cases (18); Belice moves to the front, Santa
# Lucia settles in the middle, Nueva Caledonia drops to last.
$ws.Range("A195").Value = "Belice"
$ws.Range("B195").Value = 18
$ws.Range("C195").Value = 0
$ws.Range("D195").Value = 16
$ws.Range("E195").Value = 0
$ws.Range("F195").Value = 0
$ws.Range("G195").Value = 0
$ws.Range("H195").Value = 2

$ws.Range("A196").Value = "Nueva Caledonia"
$ws.Range("B196").Value = 18
$ws.Range("C196").Value = 0
$ws.Range("D196").Value = 18
$ws.Range("E196").Value = 0
$ws.Range("F196").Value = 0
$ws.Range("G196").Value = 0
$ws.Range("H196").Value = 0

$ws.Range("A197").Value = "Santa Lucia"
$ws.Range("B197").Value = 18
$ws.Range("C197").Value = 0
$ws.Range("D197").Value = 18
$ws.Range("E197").Value = 0
$ws.Range("F197").Value = 0
$ws.Range("G197").Value = 0
$ws.Range("H197").Value = 0
